$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.204228666666666
$ws.Range("H2").Value = 24.612686
$ws.Range("I2").Value = 0.07326752815431403
$ws.Range("J2").Value = 0.07326752815431405
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.99161333333333
$ws.Range("N2").Value = 92.97484
$ws.Range("O2").Value = 0.3599121977633812
$ws.Range("P2").Value = 0.3599121977633811
$ws.Range("Q2").Value = 254.2622825355822
$ws.Range("R2").Value = 2288.36054282024
$ws.Range("S2").Value = 0.02636987708270958
$ws.Range("T2").Value = 0.02636987708270957
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.204228666666666
$ws.Range("H3").Value = 24.612686
$ws.Range("I3").Value = 0.07326752815431403
$ws.Range("J3").Value = 0.07326752815431405
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.913269
$ws.Range("N3").Value = 89.739807
$ws.Range("O3").Value = 0.3473891556493311
$ws.Range("P3").Value = 0.3473891556493311
$ws.Range("Q3").Value = 245.4152990435113
$ws.Range("R3").Value = 2208.737691391602
$ws.Range("S3").Value = 0.02545234474204075
$ws.Range("T3").Value = 0.02545234474204075
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.204228666666666
$ws.Range("H4").Value = 24.612686
$ws.Range("I4").Value = 0.07326752815431403
$ws.Range("J4").Value = 0.07326752815431405
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.150218
$ws.Range("N4").Value = 57.450654
$ws.Range("O4").Value = 0.2223955550134164
$ws.Range("P4").Value = 0.2223955550134163
$ws.Range("Q4").Value = 157.112767488516
$ws.Range("R4").Value = 1414.014907396644
$ws.Range("S4").Value = 0.01629437258833978
$ws.Range("T4").Value = 0.01629437258833978
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.204228666666666
$ws.Range("H5").Value = 24.612686
$ws.Range("I5").Value = 0.07326752815431403
$ws.Range("J5").Value = 0.07326752815431405
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.053716000000001
$ws.Range("N5").Value = 18.161148
$ws.Range("O5").Value = 0.07030309157387134
$ws.Range("P5").Value = 0.07030309157387132
$ws.Range("Q5").Value = 49.66607034705866
$ws.Range("R5").Value = 446.9946331235279
$ws.Range("S5").Value = 0.005150933741223935
$ws.Range("T5").Value = 0.005150933741223935
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 81.515531
$ws.Range("H6").Value = 244.546593
$ws.Range("I6").Value = 0.7279711116319884
$ws.Range("J6").Value = 0.7279711116319885
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.99161333333333
$ws.Range("N6").Value = 92.97484
$ws.Range("O6").Value = 0.3599121977633812
$ws.Range("P6").Value = 0.3599121977633811
$ws.Range("Q6").Value = 2526.297817413346
$ws.Range("R6").Value = 22736.68035672012
$ws.Range("S6").Value = 0.2620056826957207
$ws.Range("T6").Value = 0.2620056826957206
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 81.515531
$ws.Range("H7").Value = 244.546593
$ws.Range("I7").Value = 0.7279711116319884
$ws.Range("J7").Value = 0.7279711116319885
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.913269
$ws.Range("N7").Value = 89.739807
$ws.Range("O7").Value = 0.3473891556493311
$ws.Range("P7").Value = 0.3473891556493311
$ws.Range("Q7").Value = 2438.396006480839
$ws.Range("R7").Value = 21945.56405832755
$ws.Range("S7").Value = 0.2528892698069414
$ws.Range("T7").Value = 0.2528892698069414
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 81.515531
$ws.Range("H8").Value = 244.546593
$ws.Range("I8").Value = 0.7279711116319884
$ws.Range("J8").Value = 0.7279711116319885
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 19.150218
$ws.Range("N8").Value = 57.450654
$ws.Range("O8").Value = 0.2223955550134164
$ws.Range("P8").Value = 0.2223955550134163
$ws.Range("Q8").Value = 1561.040189035758
$ws.Range("R8").Value = 14049.36170132182
$ws.Range("S8").Value = 0.1618975394051298
$ws.Range("T8").Value = 0.1618975394051297
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 81.515531
$ws.Range("H9").Value = 244.546593
$ws.Range("I9").Value = 0.7279711116319884
$ws.Range("J9").Value = 0.7279711116319885
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.053716000000001
$ws.Range("N9").Value = 18.161148
$ws.Range("O9").Value = 0.07030309157387134
$ws.Range("P9").Value = 0.07030309157387132
$ws.Range("Q9").Value = 493.471874263196
$ws.Range("R9").Value = 4441.246868368764
$ws.Range("S9").Value = 0.05117861972419659
$ws.Range("T9").Value = 0.05117861972419659
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.363806
$ws.Range("H10").Value = 7.091418
$ws.Range("I10").Value = 0.02110987268797113
$ws.Range("J10").Value = 0.02110987268797113
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.99161333333333
$ws.Range("N10").Value = 92.97484
$ws.Range("O10").Value = 0.3599121977633812
$ws.Range("P10").Value = 0.3599121977633811
$ws.Range("Q10").Value = 73.25816154701333
$ws.Range("R10").Value = 659.32345392312
$ws.Range("S10").Value = 0.007597700673632865
$ws.Range("T10").Value = 0.007597700673632865
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.363806
$ws.Range("H11").Value = 7.091418
$ws.Range("I11").Value = 0.02110987268797113
$ws.Range("J11").Value = 0.02110987268797113
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 29.913269
$ws.Range("N11").Value = 89.739807
$ws.Range("O11").Value = 0.3473891556493311
$ws.Range("P11").Value = 0.3473891556493311
$ws.Range("Q11").Value = 70.70916474181399
$ws.Range("R11").Value = 636.382482676326
$ws.Range("S11").Value = 0.007333340848939167
$ws.Range("T11").Value = 0.007333340848939167
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.363806
$ws.Range("H12").Value = 7.091418
$ws.Range("I12").Value = 0.02110987268797113
$ws.Range("J12").Value = 0.02110987268797113
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 19.150218
$ws.Range("N12").Value = 57.450654
$ws.Range("O12").Value = 0.2223955550134164
$ws.Range("P12").Value = 0.2223955550134163
$ws.Range("Q12").Value = 45.26740020970799
$ws.Range("R12").Value = 407.406601887372
$ws.Range("S12").Value = 0.004694741852703899
$ws.Range("T12").Value = 0.004694741852703899
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.363806
$ws.Range("H13").Value = 7.091418
$ws.Range("I13").Value = 0.02110987268797113
$ws.Range("J13").Value = 0.02110987268797113
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.053716000000001
$ws.Range("N13").Value = 18.161148
$ws.Range("O13").Value = 0.07030309157387134
$ws.Range("P13").Value = 0.07030309157387132
$ws.Range("Q13").Value = 14.309810203096
$ws.Range("R13").Value = 128.788291827864
$ws.Range("S13").Value = 0.0014840893126952
$ws.Range("T13").Value = 0.0014840893126952
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 19.892761
$ws.Range("H14").Value = 59.67828300000001
$ws.Range("I14").Value = 0.1776514875257265
$ws.Range("J14").Value = 0.1776514875257265
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.99161333333333
$ws.Range("N14").Value = 92.97484
$ws.Range("O14").Value = 0.3599121977633812
$ws.Range("P14").Value = 0.3599121977633811
$ws.Range("Q14").Value = 616.5087570444134
$ws.Range("R14").Value = 5548.57881339972
$ws.Range("S14").Value = 0.06393893731131811
$ws.Range("T14").Value = 0.0639389373113181
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 19.892761
$ws.Range("H15").Value = 59.67828300000001
$ws.Range("I15").Value = 0.1776514875257265
$ws.Range("J15").Value = 0.1776514875257265
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 29.913269
$ws.Range("N15").Value = 89.739807
$ws.Range("O15").Value = 0.3473891556493311
$ws.Range("P15").Value = 0.3473891556493311
$ws.Range("Q15").Value = 595.0575109457091
$ws.Range("R15").Value = 5355.517598511382
$ws.Range("S15").Value = 0.0617142002514098
$ws.Range("T15").Value = 0.06171420025140979
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 19.892761
$ws.Range("H16").Value = 59.67828300000001
$ws.Range("I16").Value = 0.1776514875257265
$ws.Range("J16").Value = 0.1776514875257265
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 19.150218
$ws.Range("N16").Value = 57.450654
$ws.Range("O16").Value = 0.2223955550134164
$ws.Range("P16").Value = 0.2223955550134163
$ws.Range("Q16").Value = 380.9507097718981
$ws.Range("R16").Value = 3428.556387947082
$ws.Range("S16").Value = 0.03950890116724295
$ws.Range("T16").Value = 0.03950890116724295
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 19.892761
$ws.Range("H17").Value = 59.67828300000001
$ws.Range("I17").Value = 0.1776514875257265
$ws.Range("J17").Value = 0.1776514875257265
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.053716000000001
$ws.Range("N17").Value = 18.161148
$ws.Range("O17").Value = 0.07030309157387134
$ws.Range("P17").Value = 0.07030309157387132
$ws.Range("Q17").Value = 120.425125549876
$ws.Range("R17").Value = 1083.826129948884
$ws.Range("S17").Value = 0.01248944879575561
$ws.Range("T17").Value = 0.01248944879575561
